# Update product description values in column C (sharedStrings normalized: accented
# Hungarian characters replaced with unaccented ASCII equivalents) and update the
# active selection on the data_template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_template")

$ws.Range("C2:C4").Value = 'Medve dobozos natur 140g'
$ws.Range("C5:C7").Value = 'Medve dobozos szalamis 140g'
$ws.Range("C8:C10").Value = 'Medve dobozos sonkas 140g'
$ws.Range("C11:C13").Value = 'Medve dobozos zoldfuszseres 140g'
$ws.Range("C14:C16").Value = 'Medve dobozos csipospaprikas 140g'
$ws.Range("C17:C19").Value = 'Medve dobozos kolbaszos 140g'
$ws.Range("C20:C22").Value = 'Medve dobozos light 140g'
$ws.Range("C23:C25").Value = 'Medve dobozos laktozmentes 140g'
$ws.Range("C26:C28").Value = 'Medve dobozos mix 140g'
$ws.Range("C29:C31").Value = 'Medve dobozos fokhagymas-zoldfuszeres 140g'
$ws.Range("C32:C34").Value = 'Medve dobozos natur 200g'
$ws.Range("C35:C37").Value = 'Medve dobozos szalamis 200g'
$ws.Range("C38:C40").Value = 'Medve dobozos sonkas 200g'
$ws.Range("C41:C43").Value = 'Medve dobozos zoldfuszeres 200g'
$ws.Range("C44:C46").Value = 'Medve dobozos csipos 200g'
$ws.Range("C47:C49").Value = 'Medve dobozos kolbaszos 200g'
$ws.Range("C50:C52").Value = 'Medve dobozos natur 280g'
$ws.Range("C53:C55").Value = 'Medve dobozos szalamis 280g'
$ws.Range("C56:C58").Value = 'Medve dobozos sonkas 280g'
$ws.Range("C59:C61").Value = 'Medve dobozos csipos 280g'
$ws.Range("C62:C64").Value = 'Medve tomlos natur 100g'
$ws.Range("C65:C67").Value = 'Medve tomlos csipos 100g'
$ws.Range("C68:C70").Value = 'Medve tomlos zoldfuszeres 100g'
$ws.Range("C71:C73").Value = 'Medve tomlos tejszines 100g'
$ws.Range("C74:C76").Value = 'EXPORT Medve dobozos natur 140g'
$ws.Range("C77:C79").Value = 'EXPORT Medve dobozos szalamis 140g'
$ws.Range("C80:C82").Value = 'EXPORT Medve dobozos sonkas 140g'
$ws.Range("C83:C85").Value = 'EXPORT Medve dobozos zoldfuszseres 140g'
$ws.Range("C86:C88").Value = 'EXPORT Medve dobozos csipospaprikas 140g'
$ws.Range("C89:C91").Value = 'EXPORT Medve dobozos kolbaszos 140g'
$ws.Range("C92:C94").Value = 'EXPORT Medve dobozos light 140g'
$ws.Range("C95:C97").Value = 'EXPORT Medve dobozos laktozmentes 140g'
$ws.Range("C98:C100").Value = 'EXPORT Medve dobozos mix 140g'
$ws.Range("C101:C103").Value = 'EXPORT Medve dobozos fokhagymas-zoldfuszeres 140g'
$ws.Range("C104:C106").Value = 'EXPORT Medve dobozos natur 200g'
$ws.Range("C107:C109").Value = 'EXPORT Medve dobozos szalamis 200g'
$ws.Range("C110:C112").Value = 'EXPORT Medve dobozos sonkas 200g'
$ws.Range("C113:C115").Value = 'EXPORT Medve dobozos zoldfuszeres 200g'
$ws.Range("C116:C118").Value = 'EXPORT Medve dobozos csipos 200g'
$ws.Range("C119:C121").Value = 'EXPORT Medve dobozos kolbaszos 200g'
$ws.Range("C122:C124").Value = 'EXPORT Medve dobozos natur 280g'
$ws.Range("C125:C127").Value = 'EXPORT Medve dobozos szalamis 280g'
$ws.Range("C128:C130").Value = 'EXPORT Medve dobozos sonkas 280g'
$ws.Range("C131:C133").Value = 'EXPORT Medve dobozos csipos 280g'
$ws.Range("C134:C136").Value = 'EXPORT Medve tomlos natur 100g'
$ws.Range("C137:C139").Value = 'EXPORT Medve tomlos csipos 100g'
$ws.Range("C140:C142").Value = 'EXPORT Medve tomlos zoldfuszeres 100g'
$ws.Range("C143:C145").Value = 'EXPORT Medve tomlos tejszines 100g'
$ws.Range("C146:C148").Value = 'Karavan natur tomb 1,2kg'
$ws.Range("C149:C151").Value = 'Karavan fokhagymas tomb 1,2kg'
$ws.Range("C152:C154").Value = 'Karavan sonkas tomb 1,2kg'
$ws.Range("C155:C157").Value = 'Karavan nutellas tomb 1,2kg'
$ws.Range("C158:C160").Value = 'Karavan natur tomb 2,5kg'
$ws.Range("C161:C163").Value = 'Karavan fokhagymas tomb 2,5kg'
$ws.Range("C164:C166").Value = 'Karavan sonkas tomb 2,5kg'
$ws.Range("C167:C169").Value = 'Karavan natur kordobozos 140g'
$ws.Range("C170:C172").Value = 'Karavan fokhagymas kordobozos 140g'
$ws.Range("C173:C175").Value = 'Karavan sonkas kordobozos 140g'
$ws.Range("C176:C178").Value = 'Karavan natur tomlos 100g'
$ws.Range("C179:C181").Value = 'Karavan fokhagymas tomlos 100g'
$ws.Range("C182:C184").Value = 'Karavan sonkas tomlos 100g'
$ws.Range("C185:C187").Value = 'Pannonia tomb 1,2kg'
$ws.Range("C188:C190").Value = 'Pannonia light tomb 1,2kg'
$ws.Range("C191:C193").Value = 'Pannonia barsony tomb 1,2kg'

$ws.Range("C96").Select()
